{"js": "const pairs = [\n  [\"2025-05-01 Thursday\", \"2025-05-02 Friday\"],\n  [\"320\u00d73=960\", \"387\u00d76=2322\"],\n  [\"461\u00d74=1844\", \"352\u00d73=1056\"],\n  [\"201\u00d77=1407\", \"414\u00d75=2070\"],\n  [\"808\u00d74=3232\", \"710\u00d79=6390\"],\n  [\"183\u00d72=366\", \"570\u00d77=3990\"],\n  [\"890\u00d75=4450\", \"909\u00d74=3636\"],\n  [\"548\u00d73=1644\", \"595\u00d78=4760\"],\n  [\"909\u00d76=5454\", \"372\u00d74=1488\"],\n  [\"586\u00d72=1172\", \"876\u00d78=7008\"],\n  [\"507\u00d75=2535\", \"372\u00d78=2976\"],\n  [\"653\u00d77=4571\", \"140\u00d79=1260\"],\n  [\"324\u00d78=2592\", \"515\u00d73=1545\"],\n  [\"196\u00d79=1764\", \"349\u00d78=2792\"],\n  [\"853\u00d74=3412\", \"642\u00d79=5778\"],\n  [\"492\u00d77=3444\", \"160\u00d77=1120\"],\n  [\"923\u00d79=8307\", \"481\u00d76=2886\"],\n  [\"918\u00d74=3672\", \"918\u00d79=8262\"],\n  [\"631\u00d77=4417\", \"138\u00d74=552\"],\n  [\"586\u00d77=4102\", \"678\u00d76=4068\"],\n  [\"872\u00d76=5232\", \"635\u00d74=2540\"],\n  [\"985\u00d78=7880\", \"798\u00d74=3192\"],\n  [\"661\u00d79=5949\", \"307\u00d74=1228\"],\n  [\"407\u00d76=2442\", \"982\u00d79=8838\"],\n  [\"219\u00d73=657\", \"345\u00d77=2415\"],\n  [\"983\u00d72=1966\", \"388\u00d73=1164\"]\n];\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  @(\"2025-05-01 Thursday\", \"2025-05-02 Friday\"),\n  @(\"320\u00d73=960\", \"387\u00d76=2322\"),\n  @(\"461\u00d74=1844\", \"352\u00d73=1056\"),\n  @(\"201\u00d77=1407\", \"414\u00d75=2070\"),\n  @(\"808\u00d74=3232\", \"710\u00d79=6390\"),\n  @(\"183\u00d72=366\", \"570\u00d77=3990\"),\n  @(\"890\u00d75=4450\", \"909\u00d74=3636\"),\n  @(\"548\u00d73=1644\", \"595\u00d78=4760\"),\n  @(\"909\u00d76=5454\", \"372\u00d74=1488\"),\n  @(\"586\u00d72=1172\", \"876\u00d78=7008\"),\n  @(\"507\u00d75=2535\", \"372\u00d78=2976\"),\n  @(\"653\u00d77=4571\", \"140\u00d79=1260\"),\n  @(\"324\u00d78=2592\", \"515\u00d73=1545\"),\n  @(\"196\u00d79=1764\", \"349\u00d78=2792\"),\n  @(\"853\u00d74=3412\", \"642\u00d79=5778\"),\n  @(\"492\u00d77=3444\", \"160\u00d77=1120\"),\n  @(\"923\u00d79=8307\", \"481\u00d76=2886\"),\n  @(\"918\u00d74=3672\", \"918\u00d79=8262\"),\n  @(\"631\u00d77=4417\", \"138\u00d74=552\"),\n  @(\"586\u00d77=4102\", \"678\u00d76=4068\"),\n  @(\"872\u00d76=5232\", \"635\u00d74=2540\"),\n  @(\"985\u00d78=7880\", \"798\u00d74=3192\"),\n  @(\"661\u00d79=5949\", \"307\u00d74=1228\"),\n  @(\"407\u00d76=2442\", \"982\u00d79=8838\"),\n  @(\"219\u00d73=657\", \"345\u00d77=2415\"),\n  @(\"983\u00d72=1966\", \"388\u00d73=1164\")\n)\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}"}
